$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D2:D5 and K2:K5 hold plain (non-numeric) alphanumeric text - a normal
# Value assignment is stored as text (shared string) without any fuss.
$ws.Range("D2:D5").Value = "qmutxw7146636950"
$ws.Range("K2:K5").Value = "fezelm0953798973"

# AX2:AX5 need to hold the purely-numeric-looking string "5155768377"
# as TEXT (not a number), while keeping each cell's existing style
# untouched. A direct .Value/.Formula assignment of a numeric-looking
# string is auto-coerced into a number by the engine, and forcing text
# via an apostrophe/NumberFormat on the cell itself mutates its style.
# Workaround: build the literal text via a formula in a scratch cell
# that's outside the sheet's used range, copy it, and paste *values
# only* into each target cell - this carries over the text (shared
# string) type without disturbing the target cell's style/format.
$scratch = $ws.Range("BU2")
$scratch.Formula = "=""5155768377"""
$scratch.Copy()
$ws.Range("AX2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("AX3").PasteSpecial(-4163)
$ws.Range("AX4").PasteSpecial(-4163)
$ws.Range("AX5").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
